# Prelim analysis with full Eurofins data
#
# The workbook was opened and re-saved, and the data columns were widened
# to fit their (now-wider) content: column B ("Client ID" sample names like
# "Final Week 6", "SMF-E Week 6", etc.) and column P ("% Recovery" values).
# Replicate the resulting best-fit column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B -> width 28 characters (stored OOXML width == 28)
$ws.Columns.Item(2).ColumnWidth = 27.166666666666668

# Column P -> width ~9.29 characters; engine quantizes to nearest 1/6,
# 8.43 (Excel's classic default width) lands on the closest reachable value.
$ws.Columns.Item(16).ColumnWidth = 8.43
